# Nudge the four "title" textboxes (shape "CaixaDeTexto 1") on slides 8-11
# so their position/size matches the values on slide 12's equivalent shape.
# Target Left/Top are expressed in points (EMU / 12700) with enough extra
# precision to survive the single-precision round-trip used internally and
# land on the exact EMU value after conversion.

$p = $ppt.ActivePresentation

# Slide 8 - "Sumário de user evaluation"
$shape8 = $p.Slides.Item(8).Shapes.Item("CaixaDeTexto 1")
$shape8.Left  = 289.896225072442
$shape8.Width = 378.3538666677184

# Slide 9 - "Sumário de user evaluation"
$shape9 = $p.Slides.Item(9).Shapes.Item("CaixaDeTexto 1")
$shape9.Top   = 58.54574803149606
$shape9.Width = 378.2856692913386

# Slide 10 - "Resultados estáticos e destaques"
$shape10 = $p.Slides.Item(10).Shapes.Item("CaixaDeTexto 1")
$shape10.Left = 253.99055118110238
$shape10.Top  = 58.54574803149606

# Slide 11 - "Resultados estáticos e destaques"
$shape11 = $p.Slides.Item(11).Shapes.Item("CaixaDeTexto 1")
$shape11.Top  = 58.54574803149606
